# Weekly update to the Coliflor / Vega Monumental Concepción sheet:
# a new reporting date (serial 44785 = 2022-08-12) is inserted as two new
# records (Primera / Segunda quality) at the top of the historical block,
# pushing the existing rows 231:290 down to 233:292 and extending the used
# range from A1:R290 to A1:R292.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before the current row 231 (inherits formatting,
# e.g. the date-style "s=2" on column D, from the row below).
$ws.Rows.Item(231).Insert()
$ws.Rows.Item(231).Insert()

$newRows = @(
    @{ Row = 231; A = 11; B = "Vega Monumental Concepción"; C = "Bíobío"; D = 44785;
       E = 8; F = 100112008; G = "Coliflor"; H = "Sin especificar"; I = "Primera";
       J = 2000; K = 900; L = 1000; M = 950; N = "`$/unidad"; O = "Región Metropolitana";
       P = 950; Q = 1; R = "Hortaliza" },
    @{ Row = 232; A = 11; B = "Vega Monumental Concepción"; C = "Bíobío"; D = 44785;
       E = 8; F = 100112008; G = "Coliflor"; H = "Sin especificar"; I = "Segunda";
       J = 1000; K = 800; L = 800; M = 800; N = "`$/unidad"; O = "Región Metropolitana";
       P = 800; Q = 1; R = "Hortaliza" }
)

$colOrder = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R")

foreach ($rowData in $newRows) {
    $r = $rowData.Row
    for ($c = 1; $c -le $colOrder.Length; $c++) {
        $colLetter = $colOrder[$c - 1]
        $ws.Cells.Item($r, $c).Value = $rowData[$colLetter]
    }
}
